$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 21: continuation of the audit table (row 20 in "Numéro" column)
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Accessibilité"
$ws.Range("C21").Value = "menu de navigation`na mettre en page"
$ws.Range("D21").Value = "accueil et contact`nsont collés"
$ws.Range("E21").Value = "espacer les menu`npour le rendre lisible`net visible"
$ws.Range("F21").Value = "mettre une margin-left`nau menu contact"

# Match formatting used by the rest of the table
$ws.Range("A21:B21").HorizontalAlignment = -4108
$ws.Range("A21:B21").VerticalAlignment = -4108
$ws.Range("C21:F21").HorizontalAlignment = -4108
$ws.Range("C21:F21").VerticalAlignment = -4108
$ws.Range("C21:F21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 43.2

# Row 22 gets the same wrapped style carried through on E22, left empty
$ws.Range("E22").HorizontalAlignment = -4108
$ws.Range("E22").VerticalAlignment = -4108
$ws.Range("E22").WrapText = $true

# Update selection/active cell to match the new working position
$ws.Range("G21").Select()
